$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: add P1=14, Q1=15, with same style as O1 (bold/border/center)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1:O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Data rows 2-25: swap I/K and M/O values, and add P/Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
